# Replace the static "www.drpaulduenas.com" footer text with a
# configurable =website MERGEFIELD, matching the target OOXML:
#
#   <w:r><w:rPr>...</w:rPr><w:fldChar w:fldCharType="begin"/></w:r>
#   <w:r><w:rPr>...</w:rPr><w:instrText xml:space="preserve"> MERGEFIELD =website \* MERGEFORMAT </w:instrText></w:r>
#   <w:r><w:rPr>...</w:rPr><w:fldChar w:fldCharType="separate"/></w:r>
#   <w:r><w:rPr>...</w:rPr><w:t>«=website»</w:t></w:r>
#   <w:r><w:rPr>...</w:rPr><w:fldChar w:fldCharType="end"/></w:r>

$d = $word.ActiveDocument

$targetText = "www.drpaulduenas.com"

# The run formatting (rFonts Avenir Book, bold, sz 20/szCs 20) that the
# existing run already carries, and that every new run must also carry.
$runPr = '<w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'

$fieldXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>$runPr<w:fldChar w:fldCharType="begin"/></w:r>
            <w:r>$runPr<w:instrText xml:space="preserve"> MERGEFIELD =website \* MERGEFORMAT </w:instrText></w:r>
            <w:r>$runPr<w:fldChar w:fldCharType="separate"/></w:r>
            <w:r>$runPr<w:t>&#171;=website&#187;</w:t></w:r>
            <w:r>$runPr<w:fldChar w:fldCharType="end"/></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

# Search every story range (body, headers, footers, ...) for the text and
# perform the substitution wherever it is found (it lives in the default
# footer of this template).
foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers.Item($i)
        if (-not $ftr.Exists) { continue }

        $search = $ftr.Range.Duplicate
        $found = $search.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            # Insert the new field runs right after the text that is still
            # present (this keeps the paragraph's own attributes, e.g.
            # w14:paraId/rsidR and the paragraph mark run properties,
            # intact instead of Word collapsing the paragraph).
            $search.InsertXML($fieldXml) | Out-Null

            # Now remove the original run's text via Find & Replace (rather
            # than a direct Range.Text assignment) so the now-empty run is
            # dropped entirely instead of left behind as a stray <w:t></w:t>.
            $cleanup = $ftr.Range.Duplicate
            $cleanup.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $true, "", 2) | Out-Null
        }
    }
}
